# Applies updated market-price / leve-profit figures scraped by the
# scheduled runner across the ALC/ARM/BSM/CRP/CUL/LTW/WVR leve sheets.
$wb = $excel.ActiveWorkbook

# ALC!row 5 (Leve Item ID 5503)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 51.875
$ws.Range("I5").Value = 45
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 45
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 70
$ws.Range("N5").Value = -330

# ALC!row 12 (Leve Item ID 5515)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 423.0909
$ws.Range("I12").Value = 434.875
$ws.Range("J12").Value = 391.66666
$ws.Range("K12").Value = 434.875
$ws.Range("L12").Value = 391.66666
$ws.Range("M12").Value = -264.875
$ws.Range("N12").Value = -731.66666

# ALC!row 32 (Leve Item ID 5484)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 773.3333
$ws.Range("I32").Value = 537.625
$ws.Range("K32").Value = 537.625
$ws.Range("M32").Value = -211.625

# ALC!row 41 (Leve Item ID 5478)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

# ARM!row 25 (Leve Item ID 2471)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 888.8
$ws.Range("I25").Value = 532
$ws.Range("J25").Value = 1126.6666
$ws.Range("K25").Value = 532
$ws.Range("L25").Value = 1126.6666
$ws.Range("M25").Value = -130
$ws.Range("N25").Value = -1930.6666

# ARM!row 32 (Leve Item ID 44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7226.3125
$ws.Range("I32").Value = 7226.3125
$ws.Range("K32").Value = 7226.3125
$ws.Range("M32").Value = -6939.3125

# ARM!row 35 (Leve Item ID 2473)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 1774
$ws.Range("I35").Value = 1774
$ws.Range("K35").Value = 1774
$ws.Range("M35").Value = -1368

# ARM!row 45 (Leve Item ID 27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3749.5
$ws.Range("I45").Value = 3749.5
$ws.Range("K45").Value = 3749.5
$ws.Range("M45").Value = -3372.5

# ARM!row 61 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2300
$ws.Range("I61").Value = 2450
$ws.Range("K61").Value = 2450
$ws.Range("M61").Value = -2238

# ARM!row 122 (Leve Item ID 36168)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5200
$ws.Range("I122").Value = 5200
$ws.Range("K122").Value = 15600
$ws.Range("M122").Value = -13150

# ARM!row 136 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2300
$ws.Range("I136").Value = 2450
$ws.Range("K136").Value = 7350
$ws.Range("M136").Value = -4800

# BSM!row 37 (Leve Item ID 2485)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 1159.3334
$ws.Range("I37").Value = 1159.3334
$ws.Range("K37").Value = 1159.3334
$ws.Range("M37").Value = -1022.3334

# BSM!row 80 (Leve Item ID 13747)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 316.75
$ws.Range("J80").Value = 389
$ws.Range("L80").Value = 389
$ws.Range("N80").Value = -2385

# BSM!row 83 (Leve Item ID 13747)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 316.75
$ws.Range("J83").Value = 389
$ws.Range("L83").Value = 1945
$ws.Range("N83").Value = -11929

# BSM!row 94 (Leve Item ID 19939)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2350.6667
$ws.Range("I94").Value = 2220.8
$ws.Range("K94").Value = 2220.8
$ws.Range("M94").Value = -1769.8

# BSM!row 134 (Leve Item ID 43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6347.567
$ws.Range("I134").Value = 6801.154
$ws.Range("K134").Value = 20403.462
$ws.Range("M134").Value = -17868.462

# CRP!row 16 (Leve Item ID 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4872.636
$ws.Range("I16").Value = 2622.111
$ws.Range("K16").Value = 2622.111
$ws.Range("M16").Value = -2335.111

# CRP!row 99 (Leve Item ID 36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4785.4287
$ws.Range("I99").Value = 4250
$ws.Range("J99").Value = 4999.6
$ws.Range("K99").Value = 4250
$ws.Range("L99").Value = 4999.6
$ws.Range("M99").Value = -2752
$ws.Range("N99").Value = -7995.6

# CRP!row 113 (Leve Item ID 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 4872.636
$ws.Range("I113").Value = 2622.111
$ws.Range("K113").Value = 2622.111
$ws.Range("M113").Value = -452.1109999999999

# CRP!row 126 (Leve Item ID 36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4785.4287
$ws.Range("I126").Value = 4250
$ws.Range("J126").Value = 4999.6
$ws.Range("K126").Value = 12750
$ws.Range("L126").Value = 14998.8
$ws.Range("M126").Value = -10280
$ws.Range("N126").Value = -19938.8

# CUL!row 17 (Leve Item ID 4640)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 200
$ws.Range("J17").Value = 500
$ws.Range("L17").Value = 1500
$ws.Range("N17").Value = -1838

# CUL!row 23 (Leve Item ID 4858)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 290.41666
$ws.Range("J23").Value = 250.6
$ws.Range("L23").Value = 751.8
$ws.Range("N23").Value = -1221.8

# CUL!row 38 (Leve Item ID 4860)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 34.22222
$ws.Range("I38").Value = 31.714285
$ws.Range("J38").Value = 43
$ws.Range("K38").Value = 95.142855
$ws.Range("L38").Value = 129
$ws.Range("M38").Value = 251.857145
$ws.Range("N38").Value = -823

# CUL!row 86 (Leve Item ID 12892)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 435
$ws.Range("I86").Value = 298.33334
$ws.Range("J86").Value = 640
$ws.Range("K86").Value = 895.0000200000001
$ws.Range("L86").Value = 1920
$ws.Range("M86").Value = 290.9999799999999
$ws.Range("N86").Value = -4292

# CUL!row 89 (Leve Item ID 12892)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 435
$ws.Range("I89").Value = 298.33334
$ws.Range("J89").Value = 640
$ws.Range("K89").Value = 2685.00006
$ws.Range("L89").Value = 5760
$ws.Range("M89").Value = 3242.99994
$ws.Range("N89").Value = -17616

# CUL!row 98 (Leve Item ID 19843)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 291.33334
$ws.Range("J98").Value = 291.33334
$ws.Range("L98").Value = 874.0000200000001
$ws.Range("N98").Value = -3870.00002

# CUL!row 129 (Leve Item ID 36054)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3707.4
$ws.Range("J129").Value = 3897.111
$ws.Range("L129").Value = 11691.333
$ws.Range("N129").Value = -21691.333

# CUL!row 130 (Leve Item ID 36058)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 499
$ws.Range("I130").Value = 499
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 1497
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = 3523
$ws.Range("N130").ClearContents()

# LTW!row 40 (Leve Item ID 36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4459.1
$ws.Range("I40").Value = 4410.222
$ws.Range("K40").Value = 4410.222
$ws.Range("M40").Value = -4274.222

# LTW!row 68 (Leve Item ID 12563)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 775
$ws.Range("J68").Value = 950
$ws.Range("L68").Value = 950
$ws.Range("N68").Value = -2448

# LTW!row 71 (Leve Item ID 12563)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 775
$ws.Range("J71").Value = 950
$ws.Range("L71").Value = 4750
$ws.Range("N71").Value = -12238

# LTW!row 93 (Leve Item ID 19993)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1501.5
$ws.Range("I93").Value = 1501.5
$ws.Range("K93").Value = 1501.5
$ws.Range("M93").Value = -253.5

# WVR!row 62 (Leve Item ID 12589)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4875
$ws.Range("I62").Value = 1500
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 1500
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -876
$ws.Range("N62").Value = -7248

# WVR!row 65 (Leve Item ID 12589)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4875
$ws.Range("I65").Value = 1500
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 7500
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -4380
$ws.Range("N65").Value = -36240

# WVR!row 74 (Leve Item ID 19022)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 25000
$ws.Range("I74").Value = 25000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 25000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -24064
$ws.Range("N74").ClearContents()

# WVR!row 77 (Leve Item ID 19022)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 25000
$ws.Range("I77").Value = 25000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 75000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -70320
$ws.Range("N77").ClearContents()

# WVR!row 81 (Leve Item ID 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2233.3333
$ws.Range("I81").Value = 700
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 1400
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -339
$ws.Range("N81").Value = -8122

# WVR!row 84 (Leve Item ID 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2233.3333
$ws.Range("I84").Value = 700
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 7000
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -1696
$ws.Range("N84").Value = -40608
